# TC02_Canine_Filter_Breed-AmerStaffd_Neo4jData.xlsx
# "updated all canine test cases with function changes"
#
# 1) CaseDetailStat: the header row (row 1: File Name / File Type /
#    Association / Description / Format / Size) is cleared out - the
#    function that used to populate it ran with no results this time.
# 2) CaseDetailStat_Message: the Cypher query for that block turned out to
#    be empty, so a new "Cypher query should not be an empty string"
#    message row is inserted in front of the existing Neo4j connection /
#    query / output message block, and the Cypher line itself now reports
#    an empty query string instead of the MATCH ... query text.

$wb = $excel.ActiveWorkbook

# --- 1) CaseDetailStat : clear header row, but keep the (now-empty) row ---
$ws1 = $wb.Worksheets.Item("CaseDetailStat")
$ws1.Range("A1:F1").ClearContents()
# keep row 1 present (but empty) in the saved sheet data
$ws1.Rows.Item(1).OutlineLevel = 0

# --- 2) CaseDetailStat_Message : insert the new message row, then fix up
#        the Cypher line to reflect the empty query ---
$ws2 = $wb.Worksheets.Item("CaseDetailStat_Message")

# Insert a new row above row 21, shifting the last message block down
$ws2.Rows.Item(21).Insert()
$ws2.Range("A21").Value = "Cypher query should not be an empty string"

# The Cypher: line of the (now shifted) last block is row 28, followed by
# the query text row which is now empty (no query was run)
$ws2.Range("A29").Value = ""
# keep row 29 present (but empty) in the saved sheet data
$ws2.Rows.Item(29).OutlineLevel = 0
